$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.029.62'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '1.787.95'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.61'
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.545'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.24'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +3.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0682'
$ws.Range('E10').Value = '  -4.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0939'
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('D12').Value = '2.046.66'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.30'
$ws.Range('E13').Value = '  +4.06%  '
$ws.Range('D14').Value = '1.807.41'
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('D15').Value = '34.016.32'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.75'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.49'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('E20').Value = '  -1.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  -2.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.82'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.14'
$ws.Range('E26').Value = '  +1.18%  '
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('E30').Value = '  +2.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0515'
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('E32').Value = '  -0.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.57'
$ws.Range('E33').Value = '  +2.34%  '
$ws.Range('E34').Value = '  +1.87%  '
$ws.Range('D35').Value = '1.399.68'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.651'
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.35'
$ws.Range('E38').Value = '  +8.91%  '
$ws.Range('E39').Value = '  +1.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '79.97'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  +14.60%  '
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0139'
$ws.Range('E45').Value = '  +8.51%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0510'
$ws.Range('E46').Value = '  +2.82%  '
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.98'
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '107.63'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').Value = '1.948.32'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('E51').Value = '  -0.05%  '
